# GSC export refresh: roll the 90-day window forward by one day.
# - Drop the oldest date (2025-10-16) and its "Pages" count.
# - Shift every remaining date / Pages value up one row.
# - Append the new day (2026-01-14) with its Pages count (26) at the end.
#
# Dates are stored as literal text (not real Excel dates), so plain
# `.Value = "2026-01-14"` assignments must be avoided: Excel's COM layer
# auto-parses date-shaped strings into date serials and re-styles the
# cell. Shifting existing cells via Copy/PasteSpecial sidesteps that,
# since the source cells are already typed as text. For the single new
# date string that doesn't exist anywhere yet, we materialize it with a
# throwaway formula (="2026-01-14") and Copy/PasteSpecial its *value*
# into place, which also yields a plain text cell instead of a parsed date.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$lastRow = 91
$firstDataRow = 2

# Shift dates (column A) and Pages (column C) up by one row: row N gets
# what used to be in row N+1. Doing this as two bulk range copies keeps
# the original cell typing/formatting (both source ranges are already
# text/number cells, so no re-inference happens on paste).
$ws.Range("A" + ($firstDataRow + 1) + ":A" + $lastRow).Copy()
$ws.Range("A" + $firstDataRow).PasteSpecial()

$ws.Range("C" + ($firstDataRow + 1) + ":C" + $lastRow).Copy()
$ws.Range("C" + $firstDataRow).PasteSpecial()

# New last row: the newest date (as literal text) + its Pages value.
$scratch = $ws.Cells.Item(1, 5)
$scratch.Formula = "=""2026-01-14"""
$scratch.Copy()
$ws.Cells.Item($lastRow, 1).PasteSpecial(-4163)
$scratch.Clear()

$ws.Cells.Item($lastRow, 3).Value = 26
